# Update the "Pais" COVID data sheet with the latest refresh of country stats.
# Countries whose case counts grew enough to change their ranking were moved in
# the source data; here that shows up as several rows swapping which country
# name (column A) and statistics (columns B:H) they hold, plus plain numeric
# refreshes on many other rows, and a refreshed "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 16:14"

# Update country names (column A) where the shared-string order changed
$ws.Range("A36").Value = "Irak"
$ws.Range("A37").Value = "Singapur"
$ws.Range("A38").Value = "Ucrania"
$ws.Range("A80").Value = "Kenia"
$ws.Range("A81").Value = "Tayikistan"
$ws.Range("A82").Value = "El Salvador"
$ws.Range("A83").Value = "Haiti"
$ws.Range("A115").Value = "Libano"
$ws.Range("A116").Value = "Paraguay"
$ws.Range("A201").Value = "Laos"
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Dominica"
$ws.Range("A204").Value = "Fiyi"
$ws.Range("A208").Value = "Islas Malvinas"
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A212").Value = "Montserrat"
$ws.Range("A213").Value = "Seychelles"

# Update statistic values (columns B:H) for affected rows
$ws.Range("B4").Value = 2554470
$ws.Range("C4").Value = 1514
$ws.Range("D4").Value = 1068868
$ws.Range("E4").Value = 1357929
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = 127673

$ws.Range("B7").Value = 515922
$ws.Range("C7").Value = 6476
$ws.Range("D7").Value = 300648
$ws.Range("E7").Value = 199460
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 125
$ws.Range("H7").Value = 15814

$ws.Range("B17").Value = 194445
$ws.Range("C17").Value = 46
$ws.Range("D17").Value = 177500
$ws.Range("E17").Value = 7919
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 9026

$ws.Range("B30").Value = 55343
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 19143
$ws.Range("E30").Value = 35008
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 8
$ws.Range("H30").Value = 1192

$ws.Range("B36").Value = 43262
$ws.Range("C36").Value = 2069
$ws.Range("D36").Value = 19938
$ws.Range("E36").Value = 21664
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 101
$ws.Range("H36").Value = 1660

$ws.Range("B37").Value = 43246
$ws.Range("C37").Value = 291
$ws.Range("D37").Value = 36825
$ws.Range("E37").Value = 6395
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 26

$ws.Range("B38").Value = 42065
$ws.Range("C38").Value = 948
$ws.Range("D38").Value = 18701
$ws.Range("E38").Value = 22254
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 24
$ws.Range("H38").Value = 1110

$ws.Range("B39").Value = 41189
$ws.Range("C39").Value = 323
$ws.Range("D39").Value = 26864
$ws.Range("E39").Value = 12764
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 6
$ws.Range("H39").Value = 1561

$ws.Range("B62").Value = 13792
$ws.Range("C62").Value = 227
$ws.Range("D62").Value = 12338
$ws.Range("E62").Value = 1187
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 267

$ws.Range("B71").Value = 8845
$ws.Range("C71").Value = 13
$ws.Range("D71").Value = 8138
$ws.Range("E71").Value = 458
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 249

$ws.Range("B80").Value = 5811
$ws.Range("C80").Value = 278
$ws.Range("D80").Value = 1936
$ws.Range("E80").Value = 3734
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 4
$ws.Range("H80").Value = 141

$ws.Range("B81").Value = 5747
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 4331
$ws.Range("E81").Value = 1364
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 52

$ws.Range("B82").Value = 5727
$ws.Range("C82").Value = 210
$ws.Range("D82").Value = 3447
$ws.Range("E82").Value = 2137
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 10
$ws.Range("H82").Value = 143

$ws.Range("B83").Value = 5722
$ws.Range("C83").Value = 179
$ws.Range("D83").Value = 641
$ws.Range("E83").Value = 4983
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = 98

$ws.Range("B104").Value = 2330
$ws.Range("C104").Value = 5
$ws.Range("D104").Value = 2187
$ws.Range("E104").Value = 57
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 86

$ws.Range("B113").Value = 1836
$ws.Range("C113").Value = 4
$ws.Range("D113").Value = 1814
$ws.Range("E113").Value = 12
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 10

$ws.Range("B115").Value = 1719
$ws.Range("C115").Value = 22
$ws.Range("D115").Value = 1144
$ws.Range("E115").Value = 542
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 33

$ws.Range("B116").Value = 1711
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 1013
$ws.Range("E116").Value = 685
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 13

$ws.Range("B161").Value = 293
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 216
$ws.Range("E161").Value = 71
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 6

$ws.Range("B163").Value = 256
$ws.Range("C163").Value = 1
$ws.Range("D163").Value = 102
$ws.Range("E163").Value = 145
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 9

$ws.Range("B201").Value = 19
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 19
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

$ws.Range("B202").Value = 19
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 19
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("B203").Value = 18
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 18
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

$ws.Range("B204").Value = 18
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 18
$ws.Range("E204").Value = 0
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

$ws.Range("B208").Value = 13
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 13
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

$ws.Range("B209").Value = 13
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 13
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

$ws.Range("B212").Value = 11
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 10
$ws.Range("E212").Value = 0
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 1

$ws.Range("B213").Value = 11
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 11
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

